# Update minimum genome size estimates: the C-column formulas divided the
# raw base-pair count (column Z) by an assumed 30x sequencing coverage;
# the coverage assumption changed to 20x, so every "=Zn/30/1000000000"
# formula becomes "=Zn/20/1000000000" (and Excel recalculates the cached
# <v> for each cell, which a 30->20 divisor change scales by 1.5x).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,7,8,9,10,12,14,15,16,17,18,20,21,22,23,27,28,29,32,33,34,35,36)

foreach ($r in $rows) {
    $ws.Range("C$r").Formula = "=Z$r/20/1000000000"
}

# Leave the sheet's selection where the author's last edit (filling the
# C32:C36 block) would have left it.
$ws.Range("C32:C36").Select()
